# Update marksheet corrected/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: "Marking" row, Right column -> 5
$ws.Range("B11").Value = 5

# B12: "Total" row, Right column -> 70
$ws.Range("B12").Value = 70

# E12: "Total" row, Max column -> "70/140" (corrected/total marks)
$ws.Range("E12").Value = "70/140"
